# Add a new "Shagun Jhaver et al." citation paragraph to the bibliography,
# right after the "Sarah J. Jackson ... (Cambridge: MIT Press, 2020)."
# paragraph and before the "Lori Kido Lopez," paragraph. The new paragraph
# uses the same BodyText style as its neighbors, and italicizes the journal
# title ("ACM Transactions on Computer-Human Interactions").

$d = $word.ActiveDocument

# Build the new citation's plain text (curly quotes match the document's
# existing typographic convention).
$openQuote  = [char]0x201C
$closeQuote = [char]0x201D
$newParagraphText = "Shagun Jhaver, Sucheta Ghoshal, Amy Bruckman, and Eric Gilbert, " + `
    $openQuote + "Online Harassment and Content Moderation: The Case of Blocklists." + $closeQuote + `
    " ACM Transactions on Computer-Human Interactions 25, 2, Article 12 (March 2018), 33 pages. DOI: https://doi.acm.org/10.1145/3185593"

# Step 1: locate the anchor paragraph's final sentence and, in one Find &
# Replace, append a paragraph break ("^p") plus the new citation text right
# after it. Word's Find & Replace keeps the paragraph break's style
# inherited from the paragraph it splits from (BodyText here), and avoids
# leaving stray inherited character formatting behind the way
# InsertParagraphAfter()/InsertAfter() do in this host.
$anchorText = "(Cambridge: MIT Press, 2020)."
$replacement = $anchorText + "^p" + $newParagraphText

$findRange = $d.Content
$found = $findRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
if (-not $found) {
    Write-Host "ERROR: could not find anchor paragraph text '$anchorText'"
}

# Step 2: italicize just the journal title within the newly added paragraph.
$italicTarget = "ACM Transactions on Computer-Human Interactions"
$italicRange = $d.Content
$italicRange.Find.ClearFormatting()
$italicRange.Find.Replacement.ClearFormatting()
$italicRange.Find.Replacement.Font.Italic = $true
$italicized = $italicRange.Find.Execute($italicTarget, $true, $false, $false, $false, $false, $true, 1, $false, $italicTarget, 2, $true)
if (-not $italicized) {
    Write-Host "ERROR: could not find/italicize '$italicTarget'"
}

Write-Host "New citation paragraph inserted; total paragraphs: $($d.Paragraphs.Count)"
